# Update crypto price/volume data per upstream GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.248.59"
$ws.Range("E2").Value = "  -2.65%  "
$ws.Range("D3").Value = "1.933.28"
$ws.Range("E3").Value = "  -1.33%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'330.83"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").Value = "'1.005"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4739"
$ws.Range("E7").Value = "  -4.39%  "
$ws.Range("D8").Value = "'0.4066"
$ws.Range("E8").Value = "  -3.05%  "
$ws.Range("D9").Value = "'53.35"
$ws.Range("E9").Value = "  +1.19%  "
$ws.Range("E10").Value = "  -7.85%  "
$ws.Range("D11").Value = "'1.052"
$ws.Range("E11").Value = "  -3.91%  "
$ws.Range("D12").Value = "'22.35"
$ws.Range("E12").Value = "  -2.00%  "
$ws.Range("D13").Value = "1.909.91"
$ws.Range("E13").Value = "  -3.81%  "
$ws.Range("D14").Value = "'7.557"
$ws.Range("E14").Value = "  -3.44%  "
$ws.Range("D15").Value = "'6.141"
$ws.Range("E15").Value = "  -4.58%  "
$ws.Range("D16").Value = "'1.007"
$ws.Range("E16").Value = "  +0.08%  "
$ws.Range("D17").Value = "'90.29"
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("D18").Value = "'0.00001069"
$ws.Range("E18").Value = "  -2.50%  "
$ws.Range("D19").Value = "'0.06598"
$ws.Range("E19").Value = "  -1.39%  "
$ws.Range("D20").Value = "'18.27"
$ws.Range("E20").Value = "  -4.89%  "
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D22").Value = "'5.790"
$ws.Range("E22").Value = "  -2.32%  "
$ws.Range("D23").Value = "28.284.81"
$ws.Range("E23").Value = "  -2.67%  "
$ws.Range("D24").Value = "'11.45"
$ws.Range("E24").Value = "  -4.43%  "
$ws.Range("D25").Value = "'2.292"
$ws.Range("E25").Value = "  +1.42%  "
$ws.Range("D26").Value = "2.201.53"
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "'154.69"
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("D28").Value = "'20.18"
$ws.Range("E28").Value = "  -1.76%  "
$ws.Range("D29").Value = "'2.180"
$ws.Range("E29").Value = "  -2.82%  "
$ws.Range("D30").Value = "'5.795"
$ws.Range("E30").Value = "  -7.59%  "
$ws.Range("D31").Value = "'123.94"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").Value = "'0.9845"
$ws.Range("E32").Value = "  -5.31%  "
$ws.Range("D33").Value = "'0.09618"
$ws.Range("E33").Value = "  -1.97%  "
$ws.Range("D34").Value = "'1.454"
$ws.Range("E34").Value = "  -4.31%  "
$ws.Range("D35").Value = "'5.603"
$ws.Range("E35").Value = "  -3.60%  "
$ws.Range("D36").Value = "'3.648"
$ws.Range("E36").Value = "  -0.66%  "
$ws.Range("D37").Value = "'9.312"
$ws.Range("E37").Value = "  +3.79%  "
$ws.Range("D38").Value = "'0.02328"
$ws.Range("E38").Value = "  -3.84%  "
$ws.Range("D39").Value = "'0.06181"
$ws.Range("E39").Value = "  -2.53%  "
$ws.Range("D40").Value = "'1.246"
$ws.Range("E40").Value = "  -5.39%  "
$ws.Range("D41").Value = "'0.6196"
$ws.Range("E41").Value = "  -3.56%  "
$ws.Range("D42").Value = "'11.16"
$ws.Range("E42").Value = "  -1.74%  "
$ws.Range("D43").Value = "'1.006"
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("D44").Value = "'0.1911"
$ws.Range("E44").Value = "  -2.66%  "
$ws.Range("D45").Value = "'1.314"
$ws.Range("E45").Value = "  -3.79%  "
$ws.Range("D46").Value = "'0.5917"
$ws.Range("E46").Value = "  -4.49%  "
$ws.Range("E47").Value = "  -3.33%  "
$ws.Range("D48").Value = "'2.051"
$ws.Range("E48").Value = "  -6.22%  "
$ws.Range("D49").Value = "'3.465"
$ws.Range("E49").Value = "  +0.02%  "
$ws.Range("D50").Value = "'0.06797"
$ws.Range("E50").Value = "  -3.31%  "
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "'110.10"
$ws.Range("E51").Value = "  -1.88%  "
